# Weekly update: insert a new data row (for the most recent week) at the top
# of the data block (row 4, right after the two "old" historical rows 2-3),
# pushing the previously-existing rows 4-13 down to rows 5-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4; this shifts rows 4:13 down to 5:14,
# preserving all of their existing values/formatting.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's record.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44812
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7500
$ws.Range("N4").Value = "$/cuna 10 kilos"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 750
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"

# Make sure the new date cell uses the same date/time number format as the
# other rows in column D.
$ws.Range("D4").NumberFormat = $ws.Range("D5").NumberFormat
